# ----------------------------------------------------------------------
# "added bloom perf analysis"
#
# Adds a new "Bloom Optimization" perf comparison table (categories +
# before/after series) to Sheet1 and a clustered-column chart plotting
# it, repositions the existing "PBR vs. Blinn-Phong" chart further down
# the sheet to make room, and moves the active selection.
# ----------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1. New perf data table: row 12 = category labels, row 13 = "after"
#    series, row 14 = "before" series, columns A..H.
# ------------------------------------------------------------------
$categories = @("Avocado", "Fish", "Boom Box", "Corset", "Helmet", "Lantern", "Microphone", "Telephone")
$afterOptimization  = @(2.69, 2.82, 3.23, 2.75, 4.05, 2.5499999999999998, 3.03, 2.89)
$beforeOptimization = @(3.51, 3.68, 4.09, 3.59, 4.8899999999999997, 3.39, 3.91, 3.81)

for ($i = 0; $i -lt $categories.Length; $i++) {
    $col = [char](65 + $i)
    $ws.Range("$col" + "12").Value = $categories[$i]
    $ws.Range("$col" + "13").Value = $afterOptimization[$i]
    $ws.Range("$col" + "14").Value = $beforeOptimization[$i]
}

# ------------------------------------------------------------------
# 2. Move the existing "PBR vs. Blinn-Phong" chart (Chart 3) further
#    down so the new chart can take its old spot near the top.
# ------------------------------------------------------------------
$chart3Obj = $ws.ChartObjects().Item(3)
$chart3Obj.Left = 525.9375
$chart3Obj.Top = 417.3
$chart3Obj.Width = 648.8125
$chart3Obj.Height = 386.1

# ------------------------------------------------------------------
# 3. New "Bloom Optimization" clustered column chart.
# ------------------------------------------------------------------
$chart4Shape = $ws.Shapes.AddChart2(227, 51)
$chart4Obj = $chart4Shape.Chart

$chart4Obj.SetSourceData($ws.Range("A12:H14"))
$chart4Obj.ChartType = 51

$chart4Obj.HasTitle = $true
$chart4Obj.ChartTitle.Text = "Bloom Optimization"

$chart4Obj.Legend.Position = -4107

$afterSeries = $chart4Obj.SeriesCollection(1)
$afterSeries.Name = "After optimization"
$afterSeries.HasDataLabels = $true

$beforeSeries = $chart4Obj.SeriesCollection(2)
$beforeSeries.Name = "Before optimization"
$beforeSeries.HasDataLabels = $true

$valueAxis = $chart4Obj.Axes(2)
$valueAxis.HasTitle = $true
$valueAxis.AxisTitle.Text = "Frame Time (ms)"

$chart4Obj.ChartGroups(1).GapWidth = 150

$chart4Shape.Left = 499.3
$chart4Shape.Top = 58.5
$chart4Shape.Width = 694.05
$chart4Shape.Height = 402.9

# ------------------------------------------------------------------
# 4. Move the active selection, matching the author's final cursor
#    position after adding the new content.
# ------------------------------------------------------------------
$ws.Range("T1").Select()
